$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1848.6945
$ws.Range("I19").Value = 1768.4375
$ws.Range("J19").Value = 1912.9
$ws.Range("K19").Value = 1768.4375
$ws.Range("L19").Value = 1912.9
$ws.Range("M19").Value = -1593.4375
$ws.Range("N19").Value = -2262.9
$ws.Range("H62").Value = 1335300
$ws.Range("I62").Value = 4000000
$ws.Range("J62").Value = 2950
$ws.Range("K62").Value = 4000000
$ws.Range("L62").Value = 2950
$ws.Range("M62").Value = -3999376
$ws.Range("N62").Value = -4198
$ws.Range("H65").Value = 1335300
$ws.Range("I65").Value = 4000000
$ws.Range("J65").Value = 2950
$ws.Range("K65").Value = 20000000
$ws.Range("L65").Value = 14750
$ws.Range("M65").Value = -19996880
$ws.Range("N65").Value = -20990
$ws.Range("H112").Value = 3427.889
$ws.Range("J112").Value = 3511.7307
$ws.Range("L112").Value = 10535.1921
$ws.Range("N112").Value = -12751.1921
$ws.Range("H137").Value = 3276.9092
$ws.Range("I137").Value = 2178
$ws.Range("J137").Value = 4192.6665
$ws.Range("K137").Value = 6534
$ws.Range("L137").Value = 12577.9995
$ws.Range("M137").Value = -3984
$ws.Range("N137").Value = -17677.9995
$ws.Range("H138").Value = 7847.8
$ws.Range("I138").Value = 3552.6
$ws.Range("J138").Value = 9075
$ws.Range("K138").Value = 10657.8
$ws.Range("L138").Value = 27225
$ws.Range("M138").Value = -5517.799999999999
$ws.Range("N138").Value = -37505

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6042.7017
$ws.Range("I32").Value = 630.3090999999999
$ws.Range("J32").Value = 30849.5
$ws.Range("K32").Value = 630.3090999999999
$ws.Range("L32").Value = 30849.5
$ws.Range("M32").Value = -343.3090999999999
$ws.Range("N32").Value = -31423.5
$ws.Range("H45").Value = 3555.926
$ws.Range("I45").Value = 2960
$ws.Range("J45").Value = 4300.8335
$ws.Range("K45").Value = 2960
$ws.Range("L45").Value = 4300.8335
$ws.Range("M45").Value = -2583
$ws.Range("N45").Value = -5054.8335
$ws.Range("H74").Value = 2577.5
$ws.Range("I74").Value = 1962
$ws.Range("K74").Value = 1962
$ws.Range("M74").Value = -1088
$ws.Range("H77").Value = 2577.5
$ws.Range("I77").Value = 1962
$ws.Range("K77").Value = 9810
$ws.Range("M77").Value = -5442
$ws.Range("H132").Value = 4879.9653
$ws.Range("I132").Value = 4741.84
$ws.Range("J132").Value = 5743.25
$ws.Range("K132").Value = 14225.52
$ws.Range("L132").Value = 17229.75
$ws.Range("M132").Value = -11695.52
$ws.Range("N132").Value = -22289.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1504.5294
$ws.Range("I94").Value = 1411.0625
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1411.0625
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -960.0625
$ws.Range("N94").Value = -3902
$ws.Range("H107").Value = 3151
$ws.Range("I107").Value = 2556.7144
$ws.Range("J107").Value = 3844.3333
$ws.Range("K107").Value = 2556.7144
$ws.Range("L107").Value = 3844.3333
$ws.Range("M107").Value = -636.7143999999998
$ws.Range("N107").Value = -7684.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5516.6387
$ws.Range("I31").Value = 4969.3184
$ws.Range("J31").Value = 6376.7144
$ws.Range("K31").Value = 4969.3184
$ws.Range("L31").Value = 6376.7144
$ws.Range("M31").Value = -4674.3184
$ws.Range("N31").Value = -6966.7144
$ws.Range("H34").Value = 5516.6387
$ws.Range("I34").Value = 4969.3184
$ws.Range("J34").Value = 6376.7144
$ws.Range("K34").Value = 4969.3184
$ws.Range("L34").Value = 6376.7144
$ws.Range("M34").Value = -4767.3184
$ws.Range("N34").Value = -6780.7144

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1175.8125
$ws.Range("J5").Value = 1108.875
$ws.Range("L5").Value = 3326.625
$ws.Range("N5").Value = -3550.625
$ws.Range("H68").Value = 3105.054
$ws.Range("J68").Value = 3140.818
$ws.Range("L68").Value = 9422.454000000002
$ws.Range("N68").Value = -11044.454
$ws.Range("H71").Value = 3105.054
$ws.Range("J71").Value = 3140.818
$ws.Range("L71").Value = 28267.362
$ws.Range("N71").Value = -36379.362
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H113").Value = 6099.5
$ws.Range("J113").Value = 6099.5
$ws.Range("L113").Value = 18298.5
$ws.Range("N113").Value = -22638.5
$ws.Range("H130").Value = 8099.5
$ws.Range("I130").Value = 3699
$ws.Range("K130").Value = 11097
$ws.Range("M130").Value = -6077
$ws.Range("H132").Value = 4088
$ws.Range("I132").Value = 4088
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 36792
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -34262
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 1175.8125
$ws.Range("J135").Value = 1108.875
$ws.Range("L135").Value = 9979.875
$ws.Range("N135").Value = -15049.875
$ws.Range("H140").Value = 2271.25
$ws.Range("I140").Value = 2139.52
$ws.Range("J140").Value = 3369
$ws.Range("K140").Value = 6418.559999999999
$ws.Range("L140").Value = 10107
$ws.Range("M140").Value = -1238.559999999999
$ws.Range("N140").Value = -20467

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2801.9412
$ws.Range("I122").Value = 1649.5714
$ws.Range("J122").Value = 8179.6665
$ws.Range("K122").Value = 4948.7142
$ws.Range("L122").Value = 24538.9995
$ws.Range("M122").Value = -2498.7142
$ws.Range("N122").Value = -29438.9995
$ws.Range("H126").Value = 6599.8335
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 7119.8
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 21359.4
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -26299.4
$ws.Range("H132").Value = 2817.3062
$ws.Range("I132").Value = 2882.8206
$ws.Range("J132").Value = 2561.8
$ws.Range("K132").Value = 8648.461800000001
$ws.Range("L132").Value = 7685.400000000001
$ws.Range("M132").Value = -6118.461800000001
$ws.Range("N132").Value = -12745.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18031.166
$ws.Range("I7").Value = 12871.091
$ws.Range("K7").Value = 12871.091
$ws.Range("M7").Value = -12759.091
$ws.Range("H40").Value = 7006.343
$ws.Range("J40").Value = 11758.143
$ws.Range("L40").Value = 11758.143
$ws.Range("N40").Value = -12030.143
$ws.Range("H61").Value = 3997.875
$ws.Range("I61").Value = 3980.6667
$ws.Range("J61").Value = 4049.5
$ws.Range("K61").Value = 3980.6667
$ws.Range("L61").Value = 4049.5
$ws.Range("M61").Value = -3778.6667
$ws.Range("N61").Value = -4453.5
$ws.Range("H113").Value = 3997.875
$ws.Range("I113").Value = 3980.6667
$ws.Range("J113").Value = 4049.5
$ws.Range("K113").Value = 3980.6667
$ws.Range("L113").Value = 4049.5
$ws.Range("M113").Value = -1810.6667
$ws.Range("N113").Value = -8389.5
$ws.Range("H122").Value = 4496.325
$ws.Range("I122").Value = 4361.75
$ws.Range("J122").Value = 5034.625
$ws.Range("K122").Value = 13085.25
$ws.Range("L122").Value = 15103.875
$ws.Range("M122").Value = -10635.25
$ws.Range("N122").Value = -20003.875
$ws.Range("H126").Value = 18031.166
$ws.Range("I126").Value = 12871.091
$ws.Range("K126").Value = 38613.273
$ws.Range("M126").Value = -36143.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 54042.5
$ws.Range("I81").Value = 4620
$ws.Range("J81").Value = 400000
$ws.Range("K81").Value = 9240
$ws.Range("L81").Value = 800000
$ws.Range("M81").Value = -8179
$ws.Range("N81").Value = -802122
$ws.Range("H84").Value = 54042.5
$ws.Range("I84").Value = 4620
$ws.Range("J84").Value = 400000
$ws.Range("K84").Value = 46200
$ws.Range("L84").Value = 4000000
$ws.Range("M84").Value = -40896
$ws.Range("N84").Value = -4010608
$ws.Range("H125").Value = 89998.8
$ws.Range("J125").Value = 89998.8
$ws.Range("L125").Value = 89998.8
$ws.Range("N125").Value = -99838.8
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2750
$ws.Range("K126").Value = 8250
$ws.Range("M126").Value = -5780
$ws.Range("H132").Value = 2883.7878
$ws.Range("I132").Value = 2591.8965
$ws.Range("K132").Value = 7775.689499999999
$ws.Range("M132").Value = -5245.689499999999
$ws.Range("H136").Value = 3186.9756
$ws.Range("I136").Value = 2362.8125
$ws.Range("J136").Value = 6117.3335
$ws.Range("K136").Value = 7088.4375
$ws.Range("L136").Value = 18352.0005
$ws.Range("M136").Value = -4538.4375
$ws.Range("N136").Value = -23452.0005
